$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - FC Bayern München (name unchanged)
$ws.Range("B2").Value = 5.686688311688312
$ws.Range("C2").Value = 8.697986577181208
$ws.Range("D2").Value = 0.8096501572327044
$ws.Range("E2").Value = 24
$ws.Range("F2").Value = 0.04943170354428677
$ws.Range("G2").Value = 71

# Row 4 - TSG Hoffenheim (name unchanged)
$ws.Range("B4").Value = 4.584882280049566
$ws.Range("C4").Value = 8.150495049504951
$ws.Range("D4").Value = 0.5480639856480789
$ws.Range("E4").Value = 11
$ws.Range("F4").Value = 0.0842441740565065
$ws.Range("G4").Value = 10

# Row 5 - name changes from VfB Stuttgart to RB Leipzig
$ws.Range("A5").Value = "RB Leipzig"
$ws.Range("B5").Value = 5.792957746478873
$ws.Range("C5").Value = 7.952768729641694
$ws.Range("D5").Value = 0.5308977846871357
$ws.Range("E5").Value = 19
$ws.Range("F5").Value = 0.07214939168159955
$ws.Range("G5").Value = 19

# Row 6 - name changes from Bayer 04 Leverkusen to VfB Stuttgart
$ws.Range("A6").Value = "VfB Stuttgart"
$ws.Range("B6").Value = 4.919540229885057
$ws.Range("C6").Value = 8.072124756335283
$ws.Range("D6").Value = 0.6506866734486266
$ws.Range("E6").Value = 11
$ws.Range("F6").Value = 0.05817351598173516
$ws.Range("G6").Value = 11

# Row 7 - name changes from RB Leipzig to Bayer 04 Leverkusen
$ws.Range("A7").Value = "Bayer 04 Leverkusen"
$ws.Range("B7").Value = 5.888524590163934
$ws.Range("C7").Value = 10.74131274131274
$ws.Range("D7").Value = 0.6220499372933531
$ws.Range("F7").Value = 0.05697272068882273
$ws.Range("G7").Value = 26

# Row 9
$ws.Range("G9").Value = -7

# Row 10
$ws.Range("G10").Value = 0

# Row 11
$ws.Range("B11").Value = 6.455538221528861
$ws.Range("C11").Value = 8.92824074074074
$ws.Range("D11").Value = 0.4514075541094921
$ws.Range("E11").Value = 11
$ws.Range("F11").Value = 0.07709985062622084
$ws.Range("G11").Value = -5

# Row 13
$ws.Range("G13").Value = -7

# Row 14
$ws.Range("G14").Value = -27

# Row 15
$ws.Range("G15").Value = -2

# Row 18
$ws.Range("D18").Value = 0.4466342462247249
